# Fill in the Sudoku Puzzle worksheet with the solved values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final solved/updated grid state for the 9x9 Sudoku (rows 1-9, columns A-I).
# $null entries represent cells that should remain/become empty.
$grid = @(
    @(1, $null, 5, 6, 4, $null, 7, 3, 8),
    @(6, 2, 7, 1, $null, 3, 9, 4, $null),
    @(3, 4, 8, 5, 7, $null, 2, 1, 6),
    @(7, 8, 9, 4, 6, 5, 1, 2, 3),
    @(4, 6, 2, 3, 1, 7, 8, 5, $null),
    @(5, 3, 1, $null, 9, 8, 6, 7, 4),
    @(8, 1, $null, 7, 5, $null, $null, 9, 2),
    @(9, 7, 3, 8, 2, 4, 5, 6, 1),
    @(2, 5, 6, 9, 3, 1, 4, 8, 7)
)

for ($r = 0; $r -lt 9; $r++) {
    $rowVals = $grid[$r]
    for ($c = 0; $c -lt 9; $c++) {
        $val = $rowVals[$c]
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}
